# Applies the commit "excel reader and setup class" to initialSetup.xlsx
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Contacts" (was a mis-shapen 4-col/3-row sheet) -> clean
# 2-col (firstName/LastName) x 4-row (header + 3 contacts) table.
# -----------------------------------------------------------------
$contacts = $wb.Worksheets.Item("Contacts")
$contacts.UsedRange.Clear()
$contacts.Range("A1").Value = "firstName"
$contacts.Range("B1").Value = "LastName"
$contacts.Range("A2").Value = "JhordanAPI-Contact"
$contacts.Range("B2").Value = "SotoAPI-Contact"
$contacts.Range("A3").Value = "RodrigoAPI-Contact"
$contacts.Range("B3").Value = "HuancaAPI-Contact"
$contacts.Range("A4").Value = "ElyAPI-Contact"
$contacts.Range("B4").Value = "BravoAPI-Contact"
$contacts.Columns.Item(1).ColumnWidth = 21.8

# -----------------------------------------------------------------
# Sheet "Opportunities" (was a single 8-col row) -> 4-col header +
# 1 data row table.
# -----------------------------------------------------------------
$opps = $wb.Worksheets.Item("Opportunities")
$opps.UsedRange.Clear()
$opps.Range("A1").Value = "Name"
$opps.Range("B1").Value = "StageName"
$opps.Range("C1").Value = "CloseDate"
$opps.Range("D1").Value = "Amount"
$opps.Range("A2").Value = "Opportunity-createdByAPI"
$opps.Range("B2").Value = "Qualification"
$opps.Range("C2").NumberFormat = "@"
$opps.Range("C2").Value = "2022-01-18"
$opps.Range("D2").Value = 300

# -----------------------------------------------------------------
# Sheet "Campaigns" (was scratch nombre/apellido/caso data) -> a
# single Name/Value column describing the API-created campaign.
# -----------------------------------------------------------------
$campaigns = $wb.Worksheets.Item("Campaigns")
$campaigns.UsedRange.Clear()
$campaigns.Range("A1").Value = "Name"
$campaigns.Range("A2").Value = "Campaing created by API"
$campaigns.Columns.Item(1).ColumnWidth = 29.8

# -----------------------------------------------------------------
# Sheet "Cases" (was scratch nombre/apellido/caso data) -> Origin/
# Status/Priority header (highlighted font) + Phone/Working/Medium
# data row (existing "Consolas" font style).
# -----------------------------------------------------------------
$cases = $wb.Worksheets.Item("Cases")
$cases.UsedRange.Clear()
$cases.Range("A1").Value = "Origin"
$cases.Range("B1").Value = "Status"
$cases.Range("C1").Value = "Priority"
$cases.Range("A2").Value = "Phone"
$cases.Range("B2").Value = "Working"
$cases.Range("C2").Value = "Medium"

# Row 2 reuses the workbook's existing "Consolas / 9 / FF0451A5 /
# vertical-centered" style (copy the format from Opportunities!D1,
# which already carries it, so no new style entries are created).
$opps.Range("D1").Copy()
$cases.Range("A2:C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 1 needs a brand-new (but related) style: same Consolas/9/
# vertical-centered look, but colored FFA31515.
$cases.Range("A1:C1").Font.Name = "Consolas"
$cases.Range("A1:C1").Font.Size = 9
$cases.Range("A1:C1").Font.Color = 1381795
$cases.Range("A1:C1").VerticalAlignment = -4108

$cases.PageSetup.Orientation = 1

# -----------------------------------------------------------------
# Selections / active sheet: Contacts becomes the selected tab
# (workbook no longer pins activeTab=2 at Campaigns), Opportunities
# goes back to its natural A1 selection, Campaigns/Cases keep their
# own cursor positions.
# -----------------------------------------------------------------
$opps.Range("A1").Select()
$campaigns.Range("A5").Select()
$cases.Range("D6").Select()
$contacts.Range("B4").Select()
